$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 132
$ws.Range("I2").Value = 132
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 132
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -19
$ws.Range("N2").Value = ""

$ws.Range("H40").Value = 33434.285
$ws.Range("I40").Value = 27760
$ws.Range("K40").Value = 27760
$ws.Range("M40").Value = -27585

$ws.Range("H43").Value = 3856903.8
$ws.Range("J43").Value = 3000
$ws.Range("L43").Value = 3000
$ws.Range("N43").Value = -3138

$ws.Range("H70").Value = 20002608
$ws.Range("I70").Value = 1280.4
$ws.Range("J70").Value = 25002940
$ws.Range("K70").Value = 3841.2
$ws.Range("L70").Value = 75008820
$ws.Range("M70").Value = -3571.2
$ws.Range("N70").Value = -75009360

$ws.Range("H73").Value = 20002608
$ws.Range("I73").Value = 1280.4
$ws.Range("J73").Value = 25002940
$ws.Range("K73").Value = 3841.2
$ws.Range("L73").Value = 75008820
$ws.Range("M73").Value = -2905.2
$ws.Range("N73").Value = -75010692

$ws.Range("H88").Value = 15392766
$ws.Range("I88").Value = 50005124
$ws.Range("J88").Value = 9495.556
$ws.Range("K88").Value = 50005124
$ws.Range("L88").Value = 9495.556
$ws.Range("M88").Value = -50004718
$ws.Range("N88").Value = -10307.556

$ws.Range("H91").Value = 15392766
$ws.Range("I91").Value = 50005124
$ws.Range("J91").Value = 9495.556
$ws.Range("K91").Value = 50005124
$ws.Range("L91").Value = 9495.556
$ws.Range("M91").Value = -50003720
$ws.Range("N91").Value = -12303.556

$ws.Range("H97").Value = 1970
$ws.Range("J97").Value = 1970
$ws.Range("L97").Value = 5910
$ws.Range("N97").Value = -6902

$ws.Range("H100").Value = 1514.8334
$ws.Range("I100").Value = 1514.8334
$ws.Range("K100").Value = 1514.8334
$ws.Range("M100").Value = -973.8334

$ws.Range("H101").Value = 5257.6665
$ws.Range("I101").Value = 518.3333
$ws.Range("J101").Value = 9997
$ws.Range("K101").Value = 1554.9999
$ws.Range("L101").Value = 29991
$ws.Range("M101").Value = 67.00009999999997
$ws.Range("N101").Value = -33235

$ws.Range("H106").Value = 41670196
$ws.Range("I106").Value = 66668412
$ws.Range("K106").Value = 66668412
$ws.Range("M106").Value = -66667781

$ws.Range("H129").Value = 2524.1
$ws.Range("I129").Value = 726
$ws.Range("J129").Value = 5221.25
$ws.Range("K129").Value = 2178
$ws.Range("L129").Value = 15663.75
$ws.Range("M129").Value = 2822
$ws.Range("N129").Value = -25663.75

$ws.Range("H135").Value = 7008.5293
$ws.Range("I135").Value = 1124.9
$ws.Range("K135").Value = 10124.1
$ws.Range("M135").Value = -7589.1

$ws.Range("H138").Value = 5855.684
$ws.Range("J138").Value = 6581.016
$ws.Range("L138").Value = 19743.048
$ws.Range("N138").Value = -30023.048

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 730238.1
$ws.Range("J2").Value = 63770
$ws.Range("L2").Value = 63770
$ws.Range("N2").Value = -63996

$ws.Range("H32").Value = 3842.8
$ws.Range("I32").Value = 2011.0769
$ws.Range("K32").Value = 2011.0769
$ws.Range("M32").Value = -1724.0769

$ws.Range("H45").Value = 1775.2632
$ws.Range("I45").Value = 1469.0667
$ws.Range("K45").Value = 1469.0667
$ws.Range("M45").Value = -1092.0667

$ws.Range("H97").Value = 671.3125
$ws.Range("I97").Value = 696.06665
$ws.Range("K97").Value = 696.06665
$ws.Range("M97").Value = -200.06665

$ws.Range("H116").Value = 730238.1
$ws.Range("J116").Value = 63770
$ws.Range("L116").Value = 63770
$ws.Range("N116").Value = -68358

$ws.Range("H122").Value = 83338340
$ws.Range("I122").Value = 83338340
$ws.Range("K122").Value = 250015020
$ws.Range("M122").Value = -250012570

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 730238.1
$ws.Range("J3").Value = 63770
$ws.Range("L3").Value = 63770
$ws.Range("N3").Value = -63998

$ws.Range("H22").Value = 499.33334
$ws.Range("I22").Value = 499.33334
$ws.Range("K22").Value = 499.33334
$ws.Range("M22").Value = -326.33334

$ws.Range("I64").Value = 10417366
$ws.Range("J64").Value = 999
$ws.Range("K64").Value = 10417366
$ws.Range("L64").Value = 999
$ws.Range("M64").Value = -10417141
$ws.Range("N64").Value = -1449

$ws.Range("I67").Value = 10417366
$ws.Range("J67").Value = 999
$ws.Range("K67").Value = 10417366
$ws.Range("L67").Value = 999
$ws.Range("M67").Value = -10416586
$ws.Range("N67").Value = -2559

$ws.Range("H94").Value = 24462542
$ws.Range("I94").Value = 39515812
$ws.Range("K94").Value = 39515812
$ws.Range("M94").Value = -39515361

$ws.Range("H107").Value = 29413402
$ws.Range("I107").Value = 38463324
$ws.Range("K107").Value = 38463324
$ws.Range("M107").Value = -38461404

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4815.0547
$ws.Range("I31").Value = 3182.2812
$ws.Range("K31").Value = 3182.2812
$ws.Range("M31").Value = -2887.2812

$ws.Range("H34").Value = 4815.0547
$ws.Range("I34").Value = 3182.2812
$ws.Range("K34").Value = 3182.2812
$ws.Range("M34").Value = -2980.2812

$ws.Range("H62").Value = 24963.072
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376

$ws.Range("H65").Value = 24963.072
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880

$ws.Range("H99").Value = 7036.4165
$ws.Range("I99").Value = 6549.3125
$ws.Range("K99").Value = 6549.3125
$ws.Range("M99").Value = -5051.3125

$ws.Range("H126").Value = 7036.4165
$ws.Range("I126").Value = 6549.3125
$ws.Range("K126").Value = 19647.9375
$ws.Range("M126").Value = -17177.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 40563890
$ws.Range("I4").Value = 1300236
$ws.Range("J4").Value = 152745740
$ws.Range("K4").Value = 3900708
$ws.Range("L4").Value = 458237220
$ws.Range("M4").Value = -3900596
$ws.Range("N4").Value = -458237444

$ws.Range("H131").Value = 2154.35
$ws.Range("J131").Value = 2152.5208
$ws.Range("L131").Value = 6457.562399999999
$ws.Range("N131").Value = -16537.5624

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 159.5
$ws.Range("I2").Value = 177.25
$ws.Range("J2").Value = 53
$ws.Range("K2").Value = 177.25
$ws.Range("L2").Value = 53
$ws.Range("M2").Value = -64.25
$ws.Range("N2").Value = -279

$ws.Range("H70").Value = 1496124.9
$ws.Range("I70").Value = 2388418.5
$ws.Range("J70").Value = 8968.833000000001
$ws.Range("K70").Value = 2388418.5
$ws.Range("L70").Value = 8968.833000000001
$ws.Range("M70").Value = -2388148.5
$ws.Range("N70").Value = -9508.833000000001

$ws.Range("H73").Value = 1496124.9
$ws.Range("I73").Value = 2388418.5
$ws.Range("J73").Value = 8968.833000000001
$ws.Range("K73").Value = 2388418.5
$ws.Range("L73").Value = 8968.833000000001
$ws.Range("M73").Value = -2387482.5
$ws.Range("N73").Value = -10840.833

$ws.Range("H102").Value = 5369.023
$ws.Range("I102").Value = 4759.4165
$ws.Range("J102").Value = 8112.25
$ws.Range("K102").Value = 4759.4165
$ws.Range("L102").Value = 8112.25
$ws.Range("M102").Value = -3137.4165
$ws.Range("N102").Value = -11356.25

$ws.Range("H113").Value = 791627.4
$ws.Range("J113").Value = 5484.7
$ws.Range("L113").Value = 5484.7
$ws.Range("N113").Value = -9824.700000000001

$ws.Range("H122").Value = 2635358
$ws.Range("I122").Value = 4053139.8
$ws.Range("J122").Value = 2335.1428
$ws.Range("K122").Value = 12159419.4
$ws.Range("L122").Value = 7005.428400000001
$ws.Range("M122").Value = -12156969.4
$ws.Range("N122").Value = -11905.4284

$ws.Range("H126").Value = 4153.2
$ws.Range("J126").Value = 9374.75
$ws.Range("L126").Value = 28124.25
$ws.Range("N126").Value = -33064.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4935.683
$ws.Range("I7").Value = 3939.6667
$ws.Range("J7").Value = 7652.091
$ws.Range("K7").Value = 3939.6667
$ws.Range("L7").Value = 7652.091
$ws.Range("M7").Value = -3827.6667
$ws.Range("N7").Value = -7876.091

$ws.Range("H16").Value = 2510
$ws.Range("I16").Value = 1713.8
$ws.Range("J16").Value = 4500.5
$ws.Range("K16").Value = 1713.8
$ws.Range("L16").Value = 4500.5
$ws.Range("M16").Value = -1543.8
$ws.Range("N16").Value = -4840.5

$ws.Range("H126").Value = 4935.683
$ws.Range("I126").Value = 3939.6667
$ws.Range("J126").Value = 7652.091
$ws.Range("K126").Value = 11819.0001
$ws.Range("L126").Value = 22956.273
$ws.Range("M126").Value = -9349.000100000001
$ws.Range("N126").Value = -27896.273

$ws.Range("H132").Value = 5326.048
$ws.Range("I132").Value = 4359.2666
$ws.Range("K132").Value = 13077.7998
$ws.Range("M132").Value = -10547.7998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 23874.5
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").Value = ""

$ws.Range("H45").Value = 20000
$ws.Range("J45").Value = 20000
$ws.Range("L45").Value = 20000
$ws.Range("N45").Value = -20982

$ws.Range("H108").Value = 100625.4
$ws.Range("J108").Value = 100625.4
$ws.Range("L108").Value = 100625.4
$ws.Range("N108").Value = -108305.4
